# Automatische test-sync: 2025-06-17 22:13:25
# Adds a new mail log entry (row 50) to the "Logs" sheet and updates the
# "Overig" category count on the "Dashboard" sheet accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row to append after the current last row (row 49 -> row 50)
$newRow = 50

$logs.Cells.Item($newRow, 1).Value = "Sollicitatie marketingfunctie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 22:13:20"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Update the Dashboard summary count for the "Overig" category (12 -> 13)
$dashboard.Range("B3").Value = 13

# Extend the conditional formatting ranges so the new row is included too
# (Categorie column D and Beantwoord column G), mirroring D2:D49 -> D2:D50
# and G2:G49 -> G2:G50.
$dCond = $logs.Range("D2:D49").FormatConditions
$dCond.Item(1).ModifyAppliesToRange($logs.Range("D2:D50"))

$gCond = $logs.Range("G2:G49").FormatConditions
$gCond.Item(1).ModifyAppliesToRange($logs.Range("G2:G50"))
